# Insert a new worksheet "AcessarUmProdutoPelaHome_N" right after the
# existing "AcessarUmProdutoPelaHome_P" sheet, populate it with the
# "failure" variant of the product-by-category test data, and leave it as
# the active sheet/tab (matching the author's change, which added this
# sheet and a new "deveAbrirPaginaDeUmProdutoPelaCategoriaComFalha" test
# case for it).

$wb = $excel.ActiveWorkbook

$firstSheet = $wb.Worksheets.Item("AcessarUmProdutoPelaHome_P")

# Update the selection on the existing "_P" sheet (cosmetic, matches the
# committed workbook state).
$firstSheet.Select()
$firstSheet.Range("A1:C3").Select()

# Add the new sheet immediately after AcessarUmProdutoPelaHome_P so the
# final tab order becomes:
#   AcessarUmProdutoPelaHome_P, AcessarUmProdutoPelaHome_N,
#   CadastrarNovoCliente_P, CadastrarNovoCliente_N,
#   BuscarUmProdutoPelaBusca_P, BuscarUmProdutoPelaBusca_N
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $firstSheet)
$newSheet.Name = "AcessarUmProdutoPelaHome_N"

# Header row
$newSheet.Range("A1").Value = "deveAbrirPaginaDeUmProdutoPelaCategoriaComFalha"
$newSheet.Range("B1").Value = "idCategoria"
$newSheet.Range("C1").Value = "produto"

# Row 2: "headphones" category mismatched with a laptop product (failure case)
$newSheet.Range("B2").Value = "headphonesTxt"
$newSheet.Range("C2").Value = "HP Stream - 11-d020nr Laptop"

# Row 3: "laptops" category mismatched with a headphones product (failure case)
$newSheet.Range("B3").Value = "laptopsTxt"
$newSheet.Range("C3").Value = "HP H2310 In-ear Headset"

# Leave the newly created sheet selected/active, matching the committed
# workbook state (it becomes the visible tab after the edit).
$newSheet.Select()
